$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 210.33333  # H4
$ws.Cells.Item(4, 9).Value = 132.6  # I4
$ws.Cells.Item(4, 11).Value = 132.6  # K4
$ws.Cells.Item(4, 13).Value = -18.59999999999999  # M4
$ws.Cells.Item(17, 8).Value = 4010021.2  # H17
$ws.Cells.Item(32, 8).Value = 8728.368  # H32
$ws.Cells.Item(32, 9).Value = 5795  # I32
$ws.Cells.Item(32, 10).Value = 10082.23  # J32
$ws.Cells.Item(32, 11).Value = 5795  # K32
$ws.Cells.Item(32, 12).Value = 10082.23  # L32
$ws.Cells.Item(32, 13).Value = -5469  # M32
$ws.Cells.Item(32, 14).Value = -10734.23  # N32
$ws.Cells.Item(40, 8).Value = 0  # H40
$ws.Cells.Item(40, 10).Value = 0  # J40
$ws.Cells.Item(40, 12).Value = 0  # L40
$ws.Cells.Item(40, 14).ClearContents()  # N40
$ws.Cells.Item(52, 8).Value = 390  # H52
$ws.Cells.Item(52, 9).Value = 210  # I52
$ws.Cells.Item(52, 10).Value = 450  # J52
$ws.Cells.Item(52, 11).Value = 630  # K52
$ws.Cells.Item(52, 12).Value = 1350  # L52
$ws.Cells.Item(52, 13).Value = -470  # M52
$ws.Cells.Item(52, 14).Value = -1670  # N52
$ws.Cells.Item(53, 8).Value = 6496.1577  # H53
$ws.Cells.Item(53, 9).Value = 646.61536  # I53
$ws.Cells.Item(53, 11).Value = 646.61536  # K53
$ws.Cells.Item(53, 13).Value = -9.61536000000001  # M53
$ws.Cells.Item(80, 8).Value = 502.1579  # H80
$ws.Cells.Item(80, 9).Value = 571.8  # I80
$ws.Cells.Item(80, 10).Value = 424.77777  # J80
$ws.Cells.Item(80, 11).Value = 1715.4  # K80
$ws.Cells.Item(80, 12).Value = 1274.33331  # L80
$ws.Cells.Item(80, 13).Value = -717.3999999999999  # M80
$ws.Cells.Item(80, 14).Value = -3270.33331  # N80
$ws.Cells.Item(83, 8).Value = 502.1579  # H83
$ws.Cells.Item(83, 9).Value = 571.8  # I83
$ws.Cells.Item(83, 10).Value = 424.77777  # J83
$ws.Cells.Item(83, 11).Value = 5146.2  # K83
$ws.Cells.Item(83, 12).Value = 3822.99993  # L83
$ws.Cells.Item(83, 13).Value = -154.1999999999998  # M83
$ws.Cells.Item(83, 14).Value = -13806.99993  # N83
$ws.Cells.Item(86, 8).Value = 6949  # H86
$ws.Cells.Item(86, 10).Value = 9757  # J86
$ws.Cells.Item(86, 12).Value = 9757  # L86
$ws.Cells.Item(86, 14).Value = -12003  # N86
$ws.Cells.Item(88, 8).Value = 722.95654  # H88
$ws.Cells.Item(88, 10).Value = 705.3158  # J88
$ws.Cells.Item(88, 12).Value = 705.3158  # L88
$ws.Cells.Item(88, 14).Value = -1517.3158  # N88
$ws.Cells.Item(89, 8).Value = 6949  # H89
$ws.Cells.Item(89, 10).Value = 9757  # J89
$ws.Cells.Item(89, 12).Value = 48785  # L89
$ws.Cells.Item(89, 14).Value = -60017  # N89
$ws.Cells.Item(91, 8).Value = 722.95654  # H91
$ws.Cells.Item(91, 10).Value = 705.3158  # J91
$ws.Cells.Item(91, 12).Value = 705.3158  # L91
$ws.Cells.Item(91, 14).Value = -3513.3158  # N91
$ws.Cells.Item(93, 8).Value = 29333.334  # H93
$ws.Cells.Item(93, 10).Value = 29333.334  # J93
$ws.Cells.Item(93, 12).Value = 29333.334  # L93
$ws.Cells.Item(93, 14).Value = -34325.334  # N93
$ws.Cells.Item(100, 8).Value = 7360  # H100
$ws.Cells.Item(100, 9).Value = 4600  # I100
$ws.Cells.Item(100, 10).Value = 11500  # J100
$ws.Cells.Item(100, 11).Value = 4600  # K100
$ws.Cells.Item(100, 12).Value = 11500  # L100
$ws.Cells.Item(100, 13).Value = -4059  # M100
$ws.Cells.Item(100, 14).Value = -12582  # N100
$ws.Cells.Item(112, 8).Value = 6678.778  # H112
$ws.Cells.Item(112, 10).Value = 3207.8667  # J112
$ws.Cells.Item(112, 12).Value = 9623.6001  # L112
$ws.Cells.Item(112, 14).Value = -11839.6001  # N112
$ws.Cells.Item(113, 8).Value = 19402.2  # H113
$ws.Cells.Item(113, 9).Value = 18041.076  # I113
$ws.Cells.Item(113, 11).Value = 18041.076  # K113
$ws.Cells.Item(113, 13).Value = -14787.076  # M113
$ws.Cells.Item(125, 8).Value = 5372  # H125
$ws.Cells.Item(125, 9).Value = 1715  # I125
$ws.Cells.Item(125, 10).Value = 20000  # J125
$ws.Cells.Item(125, 11).Value = 15435  # K125
$ws.Cells.Item(125, 12).Value = 180000  # L125
$ws.Cells.Item(125, 13).Value = -12975  # M125
$ws.Cells.Item(125, 14).Value = -184920  # N125
$ws.Cells.Item(138, 8).Value = 5243.4736  # H138
$ws.Cells.Item(138, 10).Value = 2524.1428  # J138
$ws.Cells.Item(138, 12).Value = 7572.428400000001  # L138
$ws.Cells.Item(138, 14).Value = -17852.4284  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2993.4614  # H2
$ws.Cells.Item(2, 9).Value = 954.2632  # I2
$ws.Cells.Item(2, 10).Value = 8528.429  # J2
$ws.Cells.Item(2, 11).Value = 954.2632  # K2
$ws.Cells.Item(2, 12).Value = 8528.429  # L2
$ws.Cells.Item(2, 13).Value = -841.2632  # M2
$ws.Cells.Item(2, 14).Value = -8754.429  # N2
$ws.Cells.Item(61, 8).Value = 6302.4  # H61
$ws.Cells.Item(61, 9).Value = 1745.4359  # I61
$ws.Cells.Item(61, 11).Value = 1745.4359  # K61
$ws.Cells.Item(61, 13).Value = -1533.4359  # M61
$ws.Cells.Item(88, 8).Value = 1159  # H88
$ws.Cells.Item(88, 9).Value = 808.1667  # I88
$ws.Cells.Item(88, 11).Value = 808.1667  # K88
$ws.Cells.Item(88, 13).Value = -402.1667  # M88
$ws.Cells.Item(91, 8).Value = 1159  # H91
$ws.Cells.Item(91, 9).Value = 808.1667  # I91
$ws.Cells.Item(91, 11).Value = 808.1667  # K91
$ws.Cells.Item(91, 13).Value = 595.8333  # M91
$ws.Cells.Item(92, 8).Value = 50000  # H92
$ws.Cells.Item(92, 10).Value = 50000  # J92
$ws.Cells.Item(92, 12).Value = 50000  # L92
$ws.Cells.Item(92, 14).Value = -54992  # N92
$ws.Cells.Item(97, 8).Value = 4864.6  # H97
$ws.Cells.Item(97, 9).Value = 1497.3  # I97
$ws.Cells.Item(97, 10).Value = 11599.2  # J97
$ws.Cells.Item(97, 11).Value = 1497.3  # K97
$ws.Cells.Item(97, 12).Value = 11599.2  # L97
$ws.Cells.Item(97, 13).Value = -1001.3  # M97
$ws.Cells.Item(97, 14).Value = -12591.2  # N97
$ws.Cells.Item(102, 8).Value = 104737  # H102
$ws.Cells.Item(102, 9).Value = 9475  # I102
$ws.Cells.Item(102, 11).Value = 9475  # K102
$ws.Cells.Item(102, 13).Value = -7853  # M102
$ws.Cells.Item(110, 8).Value = 9765.666999999999  # H110
$ws.Cells.Item(110, 9).Value = 5909.4116  # I110
$ws.Cells.Item(110, 10).Value = 26154.75  # J110
$ws.Cells.Item(110, 11).Value = 5909.4116  # K110
$ws.Cells.Item(110, 12).Value = 26154.75  # L110
$ws.Cells.Item(110, 13).Value = -3864.4116  # M110
$ws.Cells.Item(110, 14).Value = -30244.75  # N110
$ws.Cells.Item(116, 8).Value = 2993.4614  # H116
$ws.Cells.Item(116, 9).Value = 954.2632  # I116
$ws.Cells.Item(116, 10).Value = 8528.429  # J116
$ws.Cells.Item(116, 11).Value = 954.2632  # K116
$ws.Cells.Item(116, 12).Value = 8528.429  # L116
$ws.Cells.Item(116, 13).Value = 1339.7368  # M116
$ws.Cells.Item(116, 14).Value = -13116.429  # N116
$ws.Cells.Item(136, 8).Value = 6302.4  # H136
$ws.Cells.Item(136, 9).Value = 1745.4359  # I136
$ws.Cells.Item(136, 11).Value = 5236.307699999999  # K136
$ws.Cells.Item(136, 13).Value = -2686.307699999999  # M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2993.4614  # H3
$ws.Cells.Item(3, 9).Value = 954.2632  # I3
$ws.Cells.Item(3, 10).Value = 8528.429  # J3
$ws.Cells.Item(3, 11).Value = 954.2632  # K3
$ws.Cells.Item(3, 12).Value = 8528.429  # L3
$ws.Cells.Item(3, 13).Value = -840.2632  # M3
$ws.Cells.Item(3, 14).Value = -8756.429  # N3
$ws.Cells.Item(94, 8).Value = 1864.4286  # H94
$ws.Cells.Item(94, 9).Value = 716.8  # I94
$ws.Cells.Item(94, 10).Value = 4733.5  # J94
$ws.Cells.Item(94, 11).Value = 716.8  # K94
$ws.Cells.Item(94, 12).Value = 4733.5  # L94
$ws.Cells.Item(94, 13).Value = -265.8  # M94
$ws.Cells.Item(94, 14).Value = -5635.5  # N94
$ws.Cells.Item(99, 8).Value = 1910.6428  # H99
$ws.Cells.Item(99, 9).Value = 1675  # I99
$ws.Cells.Item(99, 11).Value = 1675  # K99
$ws.Cells.Item(99, 13).Value = -177  # M99
$ws.Cells.Item(100, 8).Value = 31696.5  # H100
$ws.Cells.Item(100, 10).Value = 31696.5  # J100
$ws.Cells.Item(100, 12).Value = 31696.5  # L100
$ws.Cells.Item(100, 14).Value = -33860.5  # N100
$ws.Cells.Item(107, 8).Value = 1874.5883  # H107
$ws.Cells.Item(107, 9).Value = 1639.1724  # I107
$ws.Cells.Item(107, 11).Value = 1639.1724  # K107
$ws.Cells.Item(107, 13).Value = 280.8276000000001  # M107
$ws.Cells.Item(132, 8).Value = 93358.73  # H132
$ws.Cells.Item(132, 10).Value = 93358.73  # J132
$ws.Cells.Item(132, 12).Value = 93358.73  # L132
$ws.Cells.Item(132, 14).Value = -103478.73  # N132

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 4899.4  # H6
$ws.Cells.Item(6, 9).Value = 2499  # I6
$ws.Cells.Item(6, 11).Value = 2499  # K6
$ws.Cells.Item(6, 13).Value = -2386  # M6
$ws.Cells.Item(28, 8).Value = 17666.666  # H28
$ws.Cells.Item(28, 10).Value = 17666.666  # J28
$ws.Cells.Item(28, 12).Value = 17666.666  # L28
$ws.Cells.Item(28, 14).Value = -18156.666  # N28
$ws.Cells.Item(31, 8).Value = 19565.324  # H31
$ws.Cells.Item(31, 9).Value = 8993.571  # I31
$ws.Cells.Item(31, 10).Value = 26965.55  # J31
$ws.Cells.Item(31, 11).Value = 8993.571  # K31
$ws.Cells.Item(31, 12).Value = 26965.55  # L31
$ws.Cells.Item(31, 13).Value = -8698.571  # M31
$ws.Cells.Item(31, 14).Value = -27555.55  # N31
$ws.Cells.Item(34, 8).Value = 19565.324  # H34
$ws.Cells.Item(34, 9).Value = 8993.571  # I34
$ws.Cells.Item(34, 10).Value = 26965.55  # J34
$ws.Cells.Item(34, 11).Value = 8993.571  # K34
$ws.Cells.Item(34, 12).Value = 26965.55  # L34
$ws.Cells.Item(34, 13).Value = -8791.571  # M34
$ws.Cells.Item(34, 14).Value = -27369.55  # N34
$ws.Cells.Item(58, 8).Value = 11383.276  # H58
$ws.Cells.Item(58, 9).Value = 5914  # I58
$ws.Cells.Item(58, 10).Value = 13947  # J58
$ws.Cells.Item(58, 11).Value = 5914  # K58
$ws.Cells.Item(58, 12).Value = 13947  # L58
$ws.Cells.Item(58, 13).Value = -5711  # M58
$ws.Cells.Item(58, 14).Value = -14353  # N58
$ws.Cells.Item(94, 8).Value = 1174.7778  # H94
$ws.Cells.Item(94, 9).Value = 1737  # I94
$ws.Cells.Item(94, 10).Value = 893.6667  # J94
$ws.Cells.Item(94, 11).Value = 1737  # K94
$ws.Cells.Item(94, 12).Value = 893.6667  # L94
$ws.Cells.Item(94, 13).Value = -1286  # M94
$ws.Cells.Item(94, 14).Value = -1795.6667  # N94
$ws.Cells.Item(99, 8).Value = 8639.467000000001  # H99
$ws.Cells.Item(99, 10).Value = 15600.857  # J99
$ws.Cells.Item(99, 12).Value = 15600.857  # L99
$ws.Cells.Item(99, 14).Value = -18596.857  # N99
$ws.Cells.Item(105, 8).Value = 9590.706  # H105
$ws.Cells.Item(105, 9).Value = 10655.2  # I105
$ws.Cells.Item(105, 11).Value = 10655.2  # K105
$ws.Cells.Item(105, 13).Value = -8908.200000000001  # M105
$ws.Cells.Item(126, 8).Value = 8639.467000000001  # H126
$ws.Cells.Item(126, 10).Value = 15600.857  # J126
$ws.Cells.Item(126, 12).Value = 46802.571  # L126
$ws.Cells.Item(126, 14).Value = -51742.571  # N126
$ws.Cells.Item(132, 8).Value = 3887.7256  # H132
$ws.Cells.Item(132, 9).Value = 1252.3658  # I132
$ws.Cells.Item(132, 11).Value = 3757.0974  # K132
$ws.Cells.Item(132, 13).Value = -1227.0974  # M132
$ws.Cells.Item(136, 8).Value = 11383.276  # H136
$ws.Cells.Item(136, 9).Value = 5914  # I136
$ws.Cells.Item(136, 10).Value = 13947  # J136
$ws.Cells.Item(136, 11).Value = 17742  # K136
$ws.Cells.Item(136, 12).Value = 41841  # L136
$ws.Cells.Item(136, 13).Value = -15192  # M136
$ws.Cells.Item(136, 14).Value = -46941  # N136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 1248.7333  # H121
$ws.Cells.Item(121, 9).Value = 430.63635  # I121
$ws.Cells.Item(121, 10).Value = 3498.5  # J121
$ws.Cells.Item(121, 11).Value = 1291.90905  # K121
$ws.Cells.Item(121, 12).Value = 10495.5  # L121
$ws.Cells.Item(121, 13).Value = 18.09095000000002  # M121
$ws.Cells.Item(121, 14).Value = -13115.5  # N121
$ws.Cells.Item(124, 8).Value = 5997.5  # H124
$ws.Cells.Item(124, 9).Value = 5995  # I124
$ws.Cells.Item(124, 11).Value = 17985  # K124
$ws.Cells.Item(124, 13).Value = -13075  # M124
$ws.Cells.Item(132, 8).Value = 1598.7  # H132
$ws.Cells.Item(132, 10).Value = 1749  # J132
$ws.Cells.Item(132, 12).Value = 15741  # L132
$ws.Cells.Item(132, 14).Value = -20801  # N132
$ws.Cells.Item(138, 8).Value = 4054.7576  # H138
$ws.Cells.Item(138, 9).Value = 1092.1666  # I138
$ws.Cells.Item(138, 10).Value = 4713.1113  # J138
$ws.Cells.Item(138, 11).Value = 3276.4998  # K138
$ws.Cells.Item(138, 12).Value = 14139.3339  # L138
$ws.Cells.Item(138, 13).Value = 1863.5002  # M138
$ws.Cells.Item(138, 14).Value = -24419.3339  # N138
$ws.Cells.Item(140, 8).Value = 4591.25  # H140
$ws.Cells.Item(140, 9).Value = 0  # I140
$ws.Cells.Item(140, 11).Value = 0  # K140
$ws.Cells.Item(140, 13).ClearContents()  # M140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 21000  # H62
$ws.Cells.Item(62, 9).Value = 19000  # I62
$ws.Cells.Item(62, 11).Value = 19000  # K62
$ws.Cells.Item(62, 13).Value = -18314  # M62
$ws.Cells.Item(65, 8).Value = 21000  # H65
$ws.Cells.Item(65, 9).Value = 19000  # I65
$ws.Cells.Item(65, 11).Value = 57000  # K65
$ws.Cells.Item(65, 13).Value = -53568  # M65
$ws.Cells.Item(70, 8).Value = 22833  # H70
$ws.Cells.Item(70, 10).Value = 29499.5  # J70
$ws.Cells.Item(70, 12).Value = 29499.5  # L70
$ws.Cells.Item(70, 14).Value = -30039.5  # N70
$ws.Cells.Item(73, 8).Value = 22833  # H73
$ws.Cells.Item(73, 10).Value = 29499.5  # J73
$ws.Cells.Item(73, 12).Value = 29499.5  # L73
$ws.Cells.Item(73, 14).Value = -31371.5  # N73
$ws.Cells.Item(80, 8).Value = 20586.111  # H80
$ws.Cells.Item(80, 9).Value = 15476.5  # I80
$ws.Cells.Item(80, 10).Value = 24673.8  # J80
$ws.Cells.Item(80, 11).Value = 15476.5  # K80
$ws.Cells.Item(80, 12).Value = 24673.8  # L80
$ws.Cells.Item(80, 13).Value = -14478.5  # M80
$ws.Cells.Item(80, 14).Value = -26669.8  # N80
$ws.Cells.Item(83, 8).Value = 20586.111  # H83
$ws.Cells.Item(83, 9).Value = 15476.5  # I83
$ws.Cells.Item(83, 10).Value = 24673.8  # J83
$ws.Cells.Item(83, 11).Value = 77382.5  # K83
$ws.Cells.Item(83, 12).Value = 123369  # L83
$ws.Cells.Item(83, 13).Value = -72390.5  # M83
$ws.Cells.Item(83, 14).Value = -133353  # N83
$ws.Cells.Item(102, 8).Value = 10180.077  # H102
$ws.Cells.Item(102, 9).Value = 8149.222  # I102
$ws.Cells.Item(102, 11).Value = 8149.222  # K102
$ws.Cells.Item(102, 13).Value = -6527.222  # M102
$ws.Cells.Item(107, 8).Value = 865.5833  # H107
$ws.Cells.Item(107, 9).Value = 331.5  # I107
$ws.Cells.Item(107, 10).Value = 1399.6666  # J107
$ws.Cells.Item(107, 11).Value = 331.5  # K107
$ws.Cells.Item(107, 12).Value = 1399.6666  # L107
$ws.Cells.Item(107, 13).Value = 1588.5  # M107
$ws.Cells.Item(107, 14).Value = -5239.6666  # N107

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8746.723  # H7
$ws.Cells.Item(7, 10).Value = 13567.111  # J7
$ws.Cells.Item(7, 12).Value = 13567.111  # L7
$ws.Cells.Item(7, 14).Value = -13791.111  # N7
$ws.Cells.Item(40, 8).Value = 7659.05  # H40
$ws.Cells.Item(40, 9).Value = 4645.467  # I40
$ws.Cells.Item(40, 10).Value = 16699.8  # J40
$ws.Cells.Item(40, 11).Value = 4645.467  # K40
$ws.Cells.Item(40, 12).Value = 16699.8  # L40
$ws.Cells.Item(40, 13).Value = -4509.467  # M40
$ws.Cells.Item(40, 14).Value = -16971.8  # N40
$ws.Cells.Item(63, 8).Value = 20000  # H63
$ws.Cells.Item(63, 10).Value = 20500  # J63
$ws.Cells.Item(63, 12).Value = 20500  # L63
$ws.Cells.Item(63, 14).Value = -21998  # N63
$ws.Cells.Item(64, 8).Value = 77995  # H64
$ws.Cells.Item(64, 10).Value = 77995  # J64
$ws.Cells.Item(64, 12).Value = 77995  # L64
$ws.Cells.Item(64, 14).Value = -78445  # N64
$ws.Cells.Item(66, 8).Value = 20000  # H66
$ws.Cells.Item(66, 10).Value = 20500  # J66
$ws.Cells.Item(66, 12).Value = 61500  # L66
$ws.Cells.Item(66, 14).Value = -68988  # N66
$ws.Cells.Item(67, 8).Value = 77995  # H67
$ws.Cells.Item(67, 10).Value = 77995  # J67
$ws.Cells.Item(67, 12).Value = 77995  # L67
$ws.Cells.Item(67, 14).Value = -79555  # N67
$ws.Cells.Item(82, 8).Value = 4703.4136  # H82
$ws.Cells.Item(82, 10).Value = 6859.143  # J82
$ws.Cells.Item(82, 12).Value = 6859.143  # L82
$ws.Cells.Item(82, 14).Value = -7581.143  # N82
$ws.Cells.Item(85, 8).Value = 4703.4136  # H85
$ws.Cells.Item(85, 10).Value = 6859.143  # J85
$ws.Cells.Item(85, 12).Value = 6859.143  # L85
$ws.Cells.Item(85, 14).Value = -9355.143  # N85
$ws.Cells.Item(100, 8).Value = 6957.0713  # H100
$ws.Cells.Item(100, 9).Value = 3175  # I100
$ws.Cells.Item(100, 10).Value = 11999.833  # J100
$ws.Cells.Item(100, 11).Value = 3175  # K100
$ws.Cells.Item(100, 12).Value = 11999.833  # L100
$ws.Cells.Item(100, 13).Value = -2634  # M100
$ws.Cells.Item(100, 14).Value = -13081.833  # N100
$ws.Cells.Item(125, 8).Value = 129931.664  # H125
$ws.Cells.Item(125, 10).Value = 129931.664  # J125
$ws.Cells.Item(125, 12).Value = 129931.664  # L125
$ws.Cells.Item(125, 14).Value = -139771.664  # N125
$ws.Cells.Item(126, 8).Value = 8746.723  # H126
$ws.Cells.Item(126, 10).Value = 13567.111  # J126
$ws.Cells.Item(126, 12).Value = 40701.333  # L126
$ws.Cells.Item(126, 14).Value = -45641.333  # N126

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 6180.4  # H2
$ws.Cells.Item(2, 9).Value = 6180.4  # I2
$ws.Cells.Item(2, 11).Value = 6180.4  # K2
$ws.Cells.Item(2, 13).Value = -6068.4  # M2
$ws.Cells.Item(4, 8).Value = 1027.5862  # H4
$ws.Cells.Item(4, 9).Value = 1968.3636  # I4
$ws.Cells.Item(4, 10).Value = 452.66666  # J4
$ws.Cells.Item(4, 11).Value = 1968.3636  # K4
$ws.Cells.Item(4, 12).Value = 452.66666  # L4
$ws.Cells.Item(4, 13).Value = -1855.3636  # M4
$ws.Cells.Item(4, 14).Value = -678.66666  # N4
$ws.Cells.Item(75, 8).Value = 29666.334  # H75
$ws.Cells.Item(75, 9).Value = 29499.5  # I75
$ws.Cells.Item(75, 11).Value = 29499.5  # K75
$ws.Cells.Item(75, 13).Value = -28563.5  # M75
$ws.Cells.Item(78, 8).Value = 29666.334  # H78
$ws.Cells.Item(78, 9).Value = 29499.5  # I78
$ws.Cells.Item(78, 11).Value = 88498.5  # K78
$ws.Cells.Item(78, 13).Value = -83818.5  # M78
$ws.Cells.Item(81, 8).Value = 1344.6154  # H81
$ws.Cells.Item(81, 9).Value = 953.1818  # I81
$ws.Cells.Item(81, 10).Value = 3497.5  # J81
$ws.Cells.Item(81, 11).Value = 1906.3636  # K81
$ws.Cells.Item(81, 12).Value = 6995  # L81
$ws.Cells.Item(81, 13).Value = -845.3635999999999  # M81
$ws.Cells.Item(81, 14).Value = -9117  # N81
$ws.Cells.Item(84, 8).Value = 1344.6154  # H84
$ws.Cells.Item(84, 9).Value = 953.1818  # I84
$ws.Cells.Item(84, 10).Value = 3497.5  # J84
$ws.Cells.Item(84, 11).Value = 9531.817999999999  # K84
$ws.Cells.Item(84, 12).Value = 34975  # L84
$ws.Cells.Item(84, 13).Value = -4227.817999999999  # M84
$ws.Cells.Item(84, 14).Value = -45583  # N84
$ws.Cells.Item(122, 8).Value = 7222  # H122
$ws.Cells.Item(122, 9).Value = 3277  # I122
$ws.Cells.Item(122, 11).Value = 9831  # K122
$ws.Cells.Item(122, 13).Value = -7381  # M122
$ws.Cells.Item(126, 8).Value = 46781.1  # H126
$ws.Cells.Item(126, 9).Value = 58686.715  # I126
$ws.Cells.Item(126, 11).Value = 176060.145  # K126
$ws.Cells.Item(126, 13).Value = -173590.145  # M126
$ws.Cells.Item(136, 8).Value = 8471.166999999999  # H136
$ws.Cells.Item(136, 9).Value = 1898.12  # I136
$ws.Cells.Item(136, 10).Value = 23409.908  # J136
$ws.Cells.Item(136, 11).Value = 5694.36  # K136
$ws.Cells.Item(136, 12).Value = 70229.724  # L136
$ws.Cells.Item(136, 13).Value = -3144.36  # M136
$ws.Cells.Item(136, 14).Value = -75329.724  # N136
$ws.Cells.Item(138, 8).Value = 249994.5  # H138
$ws.Cells.Item(138, 10).Value = 249994.5  # J138
$ws.Cells.Item(138, 12).Value = 249994.5  # L138
$ws.Cells.Item(138, 14).Value = -260274.5  # N138
